# Filter_20211130.xlsx - 2021.11.30 完成
# Append two new stock rows (統一超 / 聯光通) to the bottom of the filter list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 80 and 81 should look exactly like the existing data rows (e.g. row 79),
# so copy the cell formatting/style from row 79 down before filling in values.
$ws.Range("A80:C81").Style = $ws.Range("A79:C79").Style

# Row 80: 統一超 (2912)
$ws.Range("A80").Value = 2912
$ws.Range("B80").Value = "統一超"
$ws.Range("C80").Value = 0

# Row 81: 聯光通 (4903)
$ws.Range("A81").Value = 4903
$ws.Range("B81").Value = "聯光通"
$ws.Range("C81").Value = 0

# Match the author's final selection/cursor position (one row below the new data).
[void]$ws.Range("A82").Select()
